# DadosEIM.xlsx - "Add files via upload" commit
# Updates the "Regramento" sheet (A1:F6) with revised wording for the
# Nivel de Servico ranges, and adjusts a few row heights + the active
# selection to match the re-saved workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Regramento")

# --- Cell content updates -------------------------------------------------
# (order matches the append order of new shared strings in the saved file)
$ws.Range("A1").Value = "Nível de Serviço"

$ws.Range("D5").Value = "D, com variação > 2 pontos percentuais"

$ws.Range("A2").Value = "[A, C]"
$ws.Range("A3").Value = "[A, C]"

$ws.Range("A4").Value = "[A, D]"
$ws.Range("B4").Value = "[A, D]"
$ws.Range("C4").Value = "[A, D]"

$ws.Range("A6").Value = "[E, H]"
$ws.Range("B6").Value = "[E, H]"
$ws.Range("C6").Value = "[E, H]"

$ws.Range("D6").Value = "[E,a H] com Variação > 1 ponto percentual"

# --- Style tweaks: B4/C4 and B6/C6 now pick up A's left border (the
# three cells visually read as one merged "range" label) ------------------
$ws.Range("A4").Copy()
$ws.Range("B4:C4").PasteSpecial(-4122)
$ws.Range("A6").Copy()
$ws.Range("B6:C6").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Row height adjustments ------------------------------------------------
$ws.Rows.Item(2).RowHeight = 90
$ws.Rows.Item(3).RowHeight = 51.75
$ws.Rows.Item(4).RowHeight = 77.25
$ws.Rows.Item(6).RowHeight = 90

# --- Selection change (matches the saved workbook's cursor position) ------
$ws.Range("D1").Select()
